# Move the four SWOT quadrant "Content Placeholder" shapes (Strengths,
# Weaknesses, Opportunities, Threats content boxes) on the slide layout
# up by 109329 EMU, keeping their horizontal position and size unchanged.
#
# Target (EMU):
#   top-row shapes    (idx 11/12 -> y 1174408 -> 1065079)
#   bottom-row shapes (idx 13/14 -> y 3892408 -> 3783079)

$p = $ppt.ActivePresentation
$master = $p.Slides.Item(1).Master
$layout = $master.CustomLayouts.Item(1)

$emuPerPoint = 12700

# A tiny epsilon nudges the COM Single-precision round-trip (pt -> EMU)
# onto the exact target EMU value instead of truncating one unit short.
$epsilon = 0.00001

$topRowPt    = (1065079 / $emuPerPoint) + $epsilon
$bottomRowPt = (3783079 / $emuPerPoint) + $epsilon

for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
    $shape = $layout.Shapes.Item($i)
    if ($shape.Name -eq "Content Placeholder 13") {
        if ([Math]::Round($shape.Top) -eq 92) {
            $shape.Top = $topRowPt
        } elseif ([Math]::Round($shape.Top) -eq 306) {
            $shape.Top = $bottomRowPt
        }
    }
}
